{"js": "// Word JS API script: update the worksheet's date header and every\n// division-problem answer cell to match the newly generated values.\n//\n// The document has one table laid out as 5 \"blocks\" of rows (a data row\n// followed by 3 blank spacer rows), 5 columns each -> 25 answer cells,\n// plus the date paragraph above the table. Every \"before\" string below\n// is unique in the document, so each edit targets exactly one run.\n\nconst body = context.document.body;\n\n// 1) Date header paragraph (e.g. \"2025-09-19 Friday\" -> \"2025-09-20 Saturday\").\nconst dateBefore = \"2025-09-19 Friday\";\nconst dateAfter = \"2025-09-20 Saturday\";\n\nconst dateResults = body.search(dateBefore, { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\n\nif (dateResults.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one match for date text '${dateBefore}', found ${dateResults.items.length}`\n  );\n}\ndateResults.items[0].insertText(dateAfter, \"Replace\");\n\n// 2) The 25 division-problem cells in the single table, addressed by\n//    (row, column) using 0-based table indices (row indices 0, 4, 8, 12,\n//    16 are the populated rows; the rows between them are blank spacers).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length !== 1) {\n  throw new Error(`Expected exactly one table, found ${tables.items.length}`);\n}\nconst table = tables.items[0];\n\n// [row, column, old cell text, new cell text]\nconst cellEdits = [\n  [0, 0, \"18\u00f79=2, 0\", \"95\u00f77=13, 4\"],\n  [0, 1, \"55\u00f76=9, 1\", \"20\u00f78=2, 4\"],\n  [0, 2, \"55\u00f72=27, 1\", \"19\u00f79=2, 1\"],\n  [0, 3, \"18\u00f74=4, 2\", \"82\u00f75=16, 2\"],\n  [0, 4, \"26\u00f73=8, 2\", \"20\u00f78=2, 4\"],\n  [4, 0, \"94\u00f74=23, 2\", \"39\u00f76=6, 3\"],\n  [4, 1, \"24\u00f75=4, 4\", \"47\u00f77=6, 5\"],\n  [4, 2, \"50\u00f79=5, 5\", \"66\u00f72=33, 0\"],\n  [4, 3, \"58\u00f79=6, 4\", \"24\u00f78=3, 0\"],\n  [4, 4, \"22\u00f78=2, 6\", \"35\u00f77=5, 0\"],\n  [8, 0, \"10\u00f77=1, 3\", \"31\u00f73=10, 1\"],\n  [8, 1, \"32\u00f72=16, 0\", \"54\u00f75=10, 4\"],\n  [8, 2, \"70\u00f77=10, 0\", \"79\u00f77=11, 2\"],\n  [8, 3, \"40\u00f74=10, 0\", \"92\u00f77=13, 1\"],\n  [8, 4, \"50\u00f72=25, 0\", \"53\u00f79=5, 8\"],\n  [12, 0, \"98\u00f79=10, 8\", \"39\u00f77=5, 4\"],\n  [12, 1, \"21\u00f77=3, 0\", \"46\u00f73=15, 1\"],\n  [12, 2, \"71\u00f77=10, 1\", \"27\u00f77=3, 6\"],\n  [12, 3, \"94\u00f78=11, 6\", \"66\u00f75=13, 1\"],\n  [12, 4, \"66\u00f78=8, 2\", \"41\u00f77=5, 6\"],\n  [16, 0, \"95\u00f75=19, 0\", \"45\u00f74=11, 1\"],\n  [16, 1, \"40\u00f76=6, 4\", \"75\u00f74=18, 3\"],\n  [16, 2, \"55\u00f75=11, 0\", \"55\u00f76=9, 1\"],\n  [16, 3, \"94\u00f79=10, 4\", \"94\u00f75=18, 4\"],\n  [16, 4, \"48\u00f79=5, 3\", \"14\u00f78=1, 6\"],\n];\n\nfor (const [row, col, before, after] of cellEdits) {\n  const cell = table.getCell(row, col);\n  cell.body.load(\"text\");\n  await context.sync();\n\n  if (cell.body.text !== before) {\n    throw new Error(\n      `Cell (${row},${col}) expected '${before}' but found '${cell.body.text}'`\n    );\n  }\n\n  cell.value = after;\n  await context.sync();\n}\n", "ps1": "# Replace the date header and all 25 division-problem answers with the\n# values from the latest generated worksheet, matching each run's old\n# text exactly before rewriting it so formatting (font/size) is kept.\n#\n# Every \"Before\" value below occurs exactly once in the document, and\n# replacements are applied in document order, so a plain Find/Replace\n# (scoped to the whole document body) on each pair in turn is safe even\n# though a couple of \"After\" values reuse text that appeared earlier\n# (that earlier occurrence has already been rewritten by the time we\n# get to it).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Before = '2025-09-19 Friday'; After = '2025-09-20 Saturday' },\n    @{ Before = '18\u00f79=2, 0'; After = '95\u00f77=13, 4' },\n    @{ Before = '55\u00f76=9, 1'; After = '20\u00f78=2, 4' },\n    @{ Before = '55\u00f72=27, 1'; After = '19\u00f79=2, 1' },\n    @{ Before = '18\u00f74=4, 2'; After = '82\u00f75=16, 2' },\n    @{ Before = '26\u00f73=8, 2'; After = '20\u00f78=2, 4' },\n    @{ Before = '94\u00f74=23, 2'; After = '39\u00f76=6, 3' },\n    @{ Before = '24\u00f75=4, 4'; After = '47\u00f77=6, 5' },\n    @{ Before = '50\u00f79=5, 5'; After = '66\u00f72=33, 0' },\n    @{ Before = '58\u00f79=6, 4'; After = '24\u00f78=3, 0' },\n    @{ Before = '22\u00f78=2, 6'; After = '35\u00f77=5, 0' },\n    @{ Before = '10\u00f77=1, 3'; After = '31\u00f73=10, 1' },\n    @{ Before = '32\u00f72=16, 0'; After = '54\u00f75=10, 4' },\n    @{ Before = '70\u00f77=10, 0'; After = '79\u00f77=11, 2' },\n    @{ Before = '40\u00f74=10, 0'; After = '92\u00f77=13, 1' },\n    @{ Before = '50\u00f72=25, 0'; After = '53\u00f79=5, 8' },\n    @{ Before = '98\u00f79=10, 8'; After = '39\u00f77=5, 4' },\n    @{ Before = '21\u00f77=3, 0'; After = '46\u00f73=15, 1' },\n    @{ Before = '71\u00f77=10, 1'; After = '27\u00f77=3, 6' },\n    @{ Before = '94\u00f78=11, 6'; After = '66\u00f75=13, 1' },\n    @{ Before = '66\u00f78=8, 2'; After = '41\u00f77=5, 6' },\n    @{ Before = '95\u00f75=19, 0'; After = '45\u00f74=11, 1' },\n    @{ Before = '40\u00f76=6, 4'; After = '75\u00f74=18, 3' },\n    @{ Before = '55\u00f75=11, 0'; After = '55\u00f76=9, 1' },\n    @{ Before = '94\u00f79=10, 4'; After = '94\u00f75=18, 4' },\n    @{ Before = '48\u00f79=5, 3'; After = '14\u00f78=1, 6' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Before\n    $find.Replacement.Text = $pair.After\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute(\n        $pair.Before,  # FindText\n        $true,         # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap = wdFindContinue\n        $false,        # Format\n        $pair.After,   # ReplaceWith\n        2              # Replace = wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"Could not find expected text: '$($pair.Before)'\"\n    }\n}\n\nWrite-Output \"Applied $($replacements.Count) replacements.\"\n"}
